$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("M94").Formula = "=[1]Treasuries!`$C`$8"
Write-Host ("M94 value=" + $ws.Range("M94").Value2)
